{"js": "// 1) Insert a new \"Meta description\" paragraph right after the title\n//    (Heading1) paragraph at the top of the document.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text,style\");\nawait context.sync();\n\nconst titlePara = paras.items[0];\n\n// Create a new, empty paragraph right after the title and reset its\n// style to Normal (it otherwise inherits the Heading1 style).\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Add the full text first, then go back and bold only the label so Word\n// ends up with two separate runs: a bold \"Meta description\" run and a\n// plain run with the rest of the sentence.\nconst metaLabel = \"Meta description\";\nconst metaBody = \": Play Book of Secrets for free and discover hidden treasures in this immersive slot game. Read our review of features, gameplay, and more.\";\nmetaPara.insertText(metaLabel + metaBody, \"End\");\nawait context.sync();\n\nconst metaLabelRange = metaPara.search(metaLabel, { matchCase: true });\nmetaLabelRange.load(\"items\");\nawait context.sync();\nmetaLabelRange.items[0].font.bold = true;\nawait context.sync();\n\n// Leave a leading empty run in the paragraph, matching the empty \"w:r\"\n// placeholder runs used throughout the rest of the document.\nmetaPara.insertText(\"\", \"Start\");\nawait context.sync();\n\n// 2) Near the end of the document, drop the paragraph that duplicated the\n//    bold title text, and rewrite the remaining (italic) paragraph so it\n//    becomes the AI image-generation prompt instead of the meta\n//    description (which now lives in the new paragraph up top).\nconst endParas = body.paragraphs;\nendParas.load(\"text\");\nawait context.sync();\n\nconst count = endParas.items.length;\nconst boldTitlePara = endParas.items[count - 2];\nconst italicPara = endParas.items[count - 1];\n\nboldTitlePara.delete();\nawait context.sync();\n\nconst oldDescription = \"Play Book of Secrets for free and discover hidden treasures in this immersive slot game. Read our review of features, gameplay, and more.\";\nconst newImagePrompt = \"Create a feature image for \\\"Book Of Secrets\\\" that features a happy Maya warrior with glasses in a cartoon style. The image should show the warrior holding the book of secrets and standing in front of a temple surrounded by jungle. The background should have vibrant colors and the warrior's expression should convey excitement and curiosity. The image should be eye-catching and make players curious about the game's mysterious theme.\";\n\nconst oldDescRange = italicPara.search(oldDescription, { matchCase: true });\noldDescRange.load(\"items\");\nawait context.sync();\noldDescRange.items[0].insertText(newImagePrompt, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Insert a new \"Meta description\" paragraph right after the title\n#    (Heading1) paragraph at the top of the document.\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Range.Style = \"Normal\"\n\n$metaLabel = \"Meta description\"\n$metaBody = \": Play Book of Secrets for free and discover hidden treasures in this immersive slot game. Read our review of features, gameplay, and more.\"\n$metaPara.Range.InsertAfter($metaLabel + $metaBody)\n\n# Bold just the \"Meta description\" label, leaving the rest of the\n# sentence as a separate, non-bold run.\n$labelRange = $metaPara.Range.Duplicate\n$labelRange.Find.Execute($metaLabel)\n$labelRange.Bold = 1\n\n# Leave a leading empty run in the paragraph, matching the empty \"w:r\"\n# placeholder runs used throughout the rest of the document.\n$startRange = $metaPara.Range.Duplicate\n$startRange.Collapse(1)\n$startRange.InsertBefore(\"\")\n\n# 2) Near the end of the document, drop the paragraph that duplicated the\n#    bold title text, and rewrite the remaining (italic) paragraph so it\n#    becomes the AI image-generation prompt instead of the meta\n#    description (which now lives in the new paragraph up top).\n$count = $d.Paragraphs.Count\n$boldTitlePara = $d.Paragraphs.Item($count - 1)\n$boldTitlePara.Range.Delete()\n\n$italicPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$oldDescription = \"Play Book of Secrets for free and discover hidden treasures in this immersive slot game. Read our review of features, gameplay, and more.\"\n$newImagePrompt = \"Create a feature image for `\"Book Of Secrets`\" that features a happy Maya warrior with glasses in a cartoon style. The image should show the warrior holding the book of secrets and standing in front of a temple surrounded by jungle. The background should have vibrant colors and the warrior's expression should convey excitement and curiosity. The image should be eye-catching and make players curious about the game's mysterious theme.\"\n\n$descRange = $italicPara.Range.Duplicate\n$descRange.Find.Execute($oldDescription)\n$descRange.Text = $newImagePrompt\n"}
